# Updated cryptos list on Sun Jul 28 17:22:54 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row, and swaps the dogwifhat/Hedera rows (44/45) to reflect their new
# ranking order. Numeric-looking Price strings are entered with a leading
# apostrophe so Excel keeps them as text (matching the workbook's existing
# inline-string cells) instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.830.72'
$ws.Range('E2').Value = '  -1.01%  '
$ws.Range('D3').Value = '3.269.38'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''579.71'
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('D6').Value = '''184.01'
$ws.Range('E6').Value = '  +1.02%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '''0.601'
$ws.Range('E8').Value = '  +0.94%  '
$ws.Range('E9').Value = '  -2.48%  '
$ws.Range('E10').Value = '  -1.12%  '
$ws.Range('D11').Value = '''0.408'
$ws.Range('E11').Value = '  -3.96%  '
$ws.Range('D12').Value = '3.837.97'
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('D14').Value = '''27.40'
$ws.Range('E14').Value = '  -4.36%  '
$ws.Range('D15').Value = '67.839.35'
$ws.Range('E15').Value = '  -1.04%  '
$ws.Range('E16').Value = '  -1.99%  '
$ws.Range('D17').Value = '3.278.29'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '''5.70'
$ws.Range('E18').Value = '  -2.30%  '
$ws.Range('D19').Value = '''13.40'
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('D20').Value = '''402.89'
$ws.Range('E20').Value = '  +2.02%  '
$ws.Range('D21').Value = '''7.54'
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('D23').Value = '''71.17'
$ws.Range('E23').Value = '  -1.16%  '
$ws.Range('E24').Value = '  -1.65%  '
$ws.Range('E25').Value = '  -1.82%  '
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('D27').Value = '''9.48'
$ws.Range('E27').Value = '  -1.80%  '
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('E29').Value = '  -1.98%  '
$ws.Range('D30').Value = '''22.70'
$ws.Range('E30').Value = '  -0.96%  '
$ws.Range('D31').Value = '''5.47'
$ws.Range('E31').Value = '  -3.96%  '
$ws.Range('D32').Value = '''6.89'
$ws.Range('E32').Value = '  -3.48%  '
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('E34').Value = '  -3.32%  '
$ws.Range('D35').Value = '''164.48'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').Value = '''1.46'
$ws.Range('E36').Value = '  -3.48%  '
$ws.Range('D37').Value = '''1.89'
$ws.Range('E37').Value = '  -1.78%  '
$ws.Range('D38').Value = '''27.07'
$ws.Range('E38').Value = '  +2.65%  '
$ws.Range('E39').Value = '  -3.13%  '
$ws.Range('D40').Value = '''4.48'
$ws.Range('E40').Value = '  -2.43%  '
$ws.Range('D41').Value = '''6.35'
$ws.Range('E41').Value = '  -2.79%  '
$ws.Range('D42').Value = '2.681.99'
$ws.Range('E42').Value = '  +2.58%  '
$ws.Range('D43').Value = '''40.88'
$ws.Range('E43').Value = '  -1.10%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').Value = '''0.0678'
$ws.Range('E44').Value = '  -1.22%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = '''2.43'
$ws.Range('E45').Value = '  -3.17%  '
$ws.Range('D46').Value = '''336.04'
$ws.Range('E46').Value = '  -2.36%  '
$ws.Range('D47').Value = '''24.58'
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('E48').Value = '  -2.79%  '
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('D51').Value = '''0.968'
$ws.Range('E51').Value = '  -1.53%  '
